# Poll Buddy Proposal (Fall 2020) - 2021-10-03 Google Drive mirror refresh.
#
# The scraped diff shows a single logical change applied uniformly across
# the whole package: every paragraph-properties block (w:pPr) -- in the
# main document body, in the header, and in the built-in heading/title
# styles -- gained an explicit <w:pageBreakBefore w:val="0"/> as its first
# child. This is exactly what toggling (or merely touching/re-saving)
# ParagraphFormat.PageBreakBefore = False does across the whole document,
# so we apply it everywhere pPr shows up: body paragraphs, header
# paragraphs, and the heading/title styles.

$d = $word.ActiveDocument

# 1) Every paragraph in the main document body.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# 2) Every paragraph in every header that already exists for each section
#    (the document has a single primary header). Only touch headers that
#    already exist so we don't mint a brand-new (empty) header/footer part
#    that wasn't in the source document.
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($p in $hdr.Range.Paragraphs) {
                $p.Format.PageBreakBefore = 0
            }
        }
    }
}

# 3) The built-in paragraph styles that carry their own pPr (the headings
#    and title/subtitle styles -- these are the ones with keepNext/keepLines
#    already set in styles.xml).
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($sn in $styleNames) {
    $s = $d.Styles.Item($sn)
    $s.ParagraphFormat.PageBreakBefore = 0
}
